## Quickstart-style refresh: rename report tabs & swap in real submission data
## (adding quickstart script for google sheets api)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: DailySiteReport (was testSheet) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "DailySiteReport"

# Header row no longer needed - clear it out completely (keeps row-level formatting)
$ws1.Range("A1:E1").Clear()

# Row 2: first real submission
$ws1.Range("A2").Value = "Robson St"
$ws1.Range("B2").Value = "Bryce Eppler"
$ws1.Range("C2").Value = "Laying pipe"
$ws1.Range("D2").Value = "12/13/1999"
$ws1.Range("E2").Value = "Today was a good day for the boys ya know haha"
$ws1.Range("A2:E2").Style = "Normal"

# Row 3: second submission
$ws1.Range("A3").Value = "DailySite"
$ws1.Range("B3").Value = "Bryce Eppler"
$ws1.Range("C3").Value = "Update"
$ws1.Range("D3").Value = "12/54/1024"
$ws1.Range("E3").Value = "comment penis"
$ws1.Range("A3:E3").Style = "Normal"

# Row 4: only the name column is populated
$ws1.Range("A4").Clear()
$ws1.Range("B4").Value = "Bryce Eppler"
$ws1.Range("B4").Style = "Normal"
$ws1.Range("C4").Clear()
$ws1.Range("D4").Clear()
$ws1.Range("E4").Clear()

# Row 5: placeholder test row
$ws1.Range("A5").Value = "test"
$ws1.Range("B5").Value = "Bryce Eppler"
$ws1.Range("C5").Value = "test"
$ws1.Range("D5").Value = "test"
$ws1.Range("E5").Value = "test"
$ws1.Range("A5:E5").Style = "Normal"

$ws1.Range("A1").Select()

# --- Sheet 2: WeeklyStaffReport (was testSheet2) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "WeeklyStaffReport"

# Drop the header labels + the lone numeric row; column G's "Date" header
# is gone too, leaving only its style behind on F1.
$ws2.Range("A1:E1").Clear()
$ws2.Range("G1").Clear()

$ws2.Range("A2").Value = "test site"
$ws2.Range("B2").Value = "Bryce Eppler"
$ws2.Range("C2").Value = "status"
$ws2.Range("D2").Value = "date"
$ws2.Range("E2").Value = "comment"
$ws2.Range("A2:E2").Style = "Normal"

$ws2.Range("A3:G3").Clear()

$ws2.Range("A1").Select()
